$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1160.2354
$ws.Range("I40").Value = 1075.8667
$ws.Range("J40").Value = 1226.8422
$ws.Range("K40").Value = 1075.8667
$ws.Range("L40").Value = 1226.8422
$ws.Range("M40").Value = -900.8667
$ws.Range("N40").Value = -1576.8422

$ws.Range("H70").Value = 8385120
$ws.Range("I70").Value = 27945850
$ws.Range("J70").Value = 1949.9286
$ws.Range("K70").Value = 83837550
$ws.Range("L70").Value = 5849.7858
$ws.Range("M70").Value = -83837280
$ws.Range("N70").Value = -6389.7858

$ws.Range("H73").Value = 8385120
$ws.Range("I73").Value = 27945850
$ws.Range("J73").Value = 1949.9286
$ws.Range("K73").Value = 83837550
$ws.Range("L73").Value = 5849.7858
$ws.Range("M73").Value = -83836614
$ws.Range("N73").Value = -7721.7858

$ws.Range("H112").Value = 4742.1665
$ws.Range("J112").Value = 4742.1665
$ws.Range("L112").Value = 14226.4995
$ws.Range("N112").Value = -16442.4995

$ws.Range("H129").Value = 1066.0667
$ws.Range("I129").Value = 307.83334
$ws.Range("J129").Value = 1150.3148
$ws.Range("K129").Value = 923.5000200000001
$ws.Range("L129").Value = 3450.9444
$ws.Range("M129").Value = 4076.49998
$ws.Range("N129").Value = -13450.9444

$ws.Range("H141").Value = 6008.913
$ws.Range("I141").Value = 1760.25
$ws.Range("K141").Value = 5280.75
$ws.Range("M141").Value = -100.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2465.476
$ws.Range("I61").Value = 2354.75
$ws.Range("K61").Value = 2354.75
$ws.Range("M61").Value = -2142.75

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H123").Value = 30429
$ws.Range("J123").Value = 30429
$ws.Range("L123").Value = 30429
$ws.Range("N123").Value = -40229

$ws.Range("H136").Value = 2465.476
$ws.Range("I136").Value = 2354.75
$ws.Range("K136").Value = 7064.25
$ws.Range("M136").Value = -4514.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1479.6111
$ws.Range("I99").Value = 1489.25
$ws.Range("J99").Value = 1460.3334
$ws.Range("K99").Value = 1489.25
$ws.Range("L99").Value = 1460.3334
$ws.Range("M99").Value = 8.75
$ws.Range("N99").Value = -4456.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 504070.2
$ws.Range("I132").Value = 714194.5600000001
$ws.Range("J132").Value = 5024.75
$ws.Range("K132").Value = 2142583.68
$ws.Range("L132").Value = 15074.25
$ws.Range("M132").Value = -2140053.68
$ws.Range("N132").Value = -20134.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 217.3
$ws.Range("I12").Value = 70.25
$ws.Range("J12").Value = 315.33334
$ws.Range("K12").Value = 210.75
$ws.Range("L12").Value = 946.0000200000001
$ws.Range("M12").Value = -37.75
$ws.Range("N12").Value = -1292.00002

$ws.Range("H107").Value = 1261.7878
$ws.Range("I107").Value = 814.9268
$ws.Range("J107").Value = 1994.64
$ws.Range("K107").Value = 2444.7804
$ws.Range("L107").Value = 5983.92
$ws.Range("M107").Value = -524.7803999999996
$ws.Range("N107").Value = -9823.92

$ws.Range("H131").Value = 2759.4036
$ws.Range("J131").Value = 3802.8975
$ws.Range("L131").Value = 11408.6925
$ws.Range("N131").Value = -21488.6925

$ws.Range("H133").Value = 7343.1665
$ws.Range("I133").Value = 2030
$ws.Range("J133").Value = 9999.75
$ws.Range("K133").Value = 6090
$ws.Range("L133").Value = 29999.25
$ws.Range("M133").Value = -1030
$ws.Range("N133").Value = -40119.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 10000
$ws.Range("I27").Value = 10000
$ws.Range("K27").Value = 10000
$ws.Range("M27").Value = -9834

$ws.Range("H62").Value = 10000
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 10000
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H70").Value = 5319.2583
$ws.Range("I70").Value = 5203.593
$ws.Range("K70").Value = 5203.593
$ws.Range("M70").Value = -4933.593

$ws.Range("H73").Value = 5319.2583
$ws.Range("I73").Value = 5203.593
$ws.Range("K73").Value = 5203.593
$ws.Range("M73").Value = -4267.593

$ws.Range("H126").Value = 3389.75
$ws.Range("I126").Value = 2517.3333
$ws.Range("K126").Value = 7551.999899999999
$ws.Range("M126").Value = -5081.999899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2471.818
$ws.Range("I68").Value = 1900
$ws.Range("J68").Value = 2948.3333
$ws.Range("K68").Value = 1900
$ws.Range("L68").Value = 2948.3333
$ws.Range("M68").Value = -1151
$ws.Range("N68").Value = -4446.3333

$ws.Range("H71").Value = 2471.818
$ws.Range("I71").Value = 1900
$ws.Range("J71").Value = 2948.3333
$ws.Range("K71").Value = 9500
$ws.Range("L71").Value = 14741.6665
$ws.Range("M71").Value = -5756
$ws.Range("N71").Value = -22229.6665

$ws.Range("H100").Value = 3055
$ws.Range("I100").Value = 2615
$ws.Range("J100").Value = 3825
$ws.Range("K100").Value = 2615
$ws.Range("L100").Value = 3825
$ws.Range("M100").Value = -2074
$ws.Range("N100").Value = -4907

$ws.Range("H122").Value = 9095670
$ws.Range("I122").Value = 4020
$ws.Range("J122").Value = 13338440
$ws.Range("K122").Value = 12060
$ws.Range("L122").Value = 40015320
$ws.Range("M122").Value = -9610
$ws.Range("N122").Value = -40020220

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5591.6
$ws.Range("I62").Value = 7000
$ws.Range("J62").Value = 4652.6665
$ws.Range("K62").Value = 7000
$ws.Range("L62").Value = 4652.6665
$ws.Range("M62").Value = -6376
$ws.Range("N62").Value = -5900.6665

$ws.Range("H65").Value = 5591.6
$ws.Range("I65").Value = 7000
$ws.Range("J65").Value = 4652.6665
$ws.Range("K65").Value = 35000
$ws.Range("L65").Value = 23263.3325
$ws.Range("M65").Value = -31880
$ws.Range("N65").Value = -29503.3325

$ws.Range("H122").Value = 1997
$ws.Range("I122").Value = 1946
$ws.Range("J122").Value = 2405
$ws.Range("K122").Value = 5838
$ws.Range("L122").Value = 7215
$ws.Range("M122").Value = -3388
$ws.Range("N122").Value = -12115

$ws.Range("H133").Value = 750000
$ws.Range("J133").Value = 750000
$ws.Range("L133").Value = 750000
$ws.Range("N133").Value = -760120

Write-Output "Asura_Profits scheduled-runner update applied"
